# "Button changed, added prep work for patient history"
#
# Observed change (from the OOXML diff):
#   1. Four new entries are appended to the shared-string table:
#        "DD.MM.YYYY", "HH.MM", "14.14", "08.10"
#      (the first three are "prep work" for an upcoming patient-history
#      feature - date/time format hints and a sample value - and are not
#      yet wired into any visible cell).
#   2. Cell E7 - previously showing "08.11" - now shows "08.10"
#      (it now points at the newly added "08.10" shared string).
#
# Because every cell in this engine is auto-type-detected, values that
# look like dates/numbers (e.g. "08.10", "14.14") must be forced to plain
# text (NumberFormat "@") before assignment, otherwise Excel would silently
# store them as numeric values instead of the text shown in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- prep work: register the new strings used by the upcoming -------------
# --- "patient history" feature so they exist in the workbook's -----------
# --- shared-string table (added in this order, matching the diff) --------
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "DD.MM.YYYY"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "HH.MM"

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "14.14"

# --- visible fix: E7 "08.11" -> "08.10" ------------------------------------
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "08.10"
